$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 656; existing rows 656:697 shift down to 657:698.
$ws.Rows.Item(656).Insert()

# Column A holds the date as literal text (e.g. "2026/12/29"), not a real
# date value, so force text formatting before writing it to avoid Excel's
# auto date-parsing, then restore the cell's style to match its sibling
# data rows (which carry no explicit style).
$ws.Range("A656").NumberFormat = "@"
$ws.Range("A656").Value = "2026/01/17"
$ws.Range("A656").Style = $ws.Range("A657").Style()

$ws.Range("B656").Value = "土"
$ws.Range("C656").Value = 17
$ws.Range("D656").Value = 201
